# feat: add 2022-Q3 data
#
#  - The "总计" (summary) sheet gets a new row for 2022-Q3 (the existing
#    2022-Q2 row shifts down to row 3).
#  - The existing "2022-Q2" detail sheet is duplicated so a pristine copy of
#    its data survives as the (new) "2022-Q2" sheet placed right after it;
#    the original sheet is renamed to "2022-Q3" and repopulated with the
#    Q3 fund-holding figures.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Duplicate the "2022-Q2" sheet so we keep a pristine copy of its data,
#    then rename the original to "2022-Q3" and the copy back to "2022-Q2".
#    (Renaming the original in place keeps its sheetId; the freshly created
#    copy picks up the next sheetId - matching how Excel numbers sheets when
#    you copy one rather than inserting a brand new blank tab.)
# ---------------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($null, $q2)
$q2copy = $wb.Worksheets.Item(3)
$q2.Name = "2022-Q3"
$q2copy.Name = "2022-Q2"

$q3 = $q2
$total = $wb.Worksheets.Item("总计")

# A cell that still carries the untouched default/General format - used below
# to strip the temporary "@" (Text) number format back off once a numeric-
# looking string has been safely stored, without disturbing the real General
# style index.
$generalFmt = $total.Range("C1")

# ---------------------------------------------------------------------------
# 2) Replace the (copied) Q2 fund-holding rows on the "2022-Q3" sheet with
#    the real Q3 figures.
# ---------------------------------------------------------------------------
$q3.Range("A1:H15").ClearContents()

$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

$q3Rows = @(
    @(0,  "001822", "华商智能生活灵活配置混合A", "33.45", "87.34", "4.51", "1.5086", 7),
    @(1,  "000729", "建信中小盘先锋股票A",       "37.13", "94.11", "2.77", "1.0285", 8),
    @(2,  "530005", "建信优化配置混合A",         "27.77", "88.12", "2.97", "0.8248", 9),
    @(3,  "001933", "华商新兴活力灵活配置混合",   "18.34", "87.96", "4.25", "0.7794", 6),
    @(4,  "013886", "华商新能源汽车混合A",       "9.72",  "89.01", "6.12", "0.5949", 6),
    @(5,  "010550", "华商双擎领航混合",           "12.53", "90.32", "4.35", "0.5451", 10),
    @(6,  "015385", "华商智能生活灵活配置混合C", "11.97", "87.34", "4.51", "0.5398", 7),
    @(7,  "000756", "建信潜力新蓝筹股票A",       "18.95", "93.65", "2.74", "0.5192", 10),
    @(8,  "013919", "建信中小盘先锋股票C",       "16.59", "94.11", "2.77", "0.4595", 8),
    @(9,  "014967", "建信潜力新蓝筹股票C",       "14.95", "93.65", "2.74", "0.4096", 10),
    @(10, "013887", "华商新能源汽车混合C",       "4.06",  "89.01", "6.12", "0.2485", 6),
    @(11, "014350", "华商卓越成长一年持有混合A", "3.14",  "86.88", "4.89", "0.1535", 6),
    @(12, "015436", "建信优化配置混合C",         "0.17",  "88.12", "2.97", "0.0050", 9),
    @(13, "014351", "华商卓越成长一年持有混合C", "0.10",  "86.88", "4.89", "0.0049", 6)
)

# Columns B and D:G hold numeric-looking text (fund codes / figures stored as
# strings in the source data). Format them as Text first so Excel doesn't
# silently coerce them into numbers (and drop leading zeros / add float
# noise), then fill in the values.
$q3.Range("B2:B15").NumberFormat = "@"
$q3.Range("D2:G15").NumberFormat = "@"

$r = 2
foreach ($row in $q3Rows) {
    $q3.Cells.Item($r, 1).Value = $row[0]
    $q3.Range("B$r").Value = $row[1]
    $q3.Range("C$r").Value = $row[2]
    $q3.Range("D$r").Value = $row[3]
    $q3.Range("E$r").Value = $row[4]
    $q3.Range("F$r").Value = $row[5]
    $q3.Range("G$r").Value = $row[6]
    $q3.Cells.Item($r, 8).Value = $row[7]
    $r++
}

# Drop the temporary Text number format back to General now that the values
# are safely stored as strings, so the cells don't carry a stray style.
$generalFmt.Copy()
$q3.Range("B2:B15").PasteSpecial(-4122)
$q3.Range("D2:G15").PasteSpecial(-4122)

# Re-apply the bold/bordered header style (copied from the "总计" sheet's
# header-row formatting) to the "2022-Q3" header row and index column.
$total.Range("A2").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)
$q3.Range("A2:A15").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3) Update the "总计" summary sheet: push the existing 2022-Q2 row down to
#    row 3 and insert the new 2022-Q3 totals in row 2.
# ---------------------------------------------------------------------------

# Move the existing 2022-Q2 totals row down to row 3 (write literals rather
# than re-reading .Value, which doesn't round-trip through this host).
$total.Range("A3").Value = 1
$total.Range("B3").Value = "2022-Q2"
$total.Range("C3").Value = 9
$total.Range("D3").Value = 2.51
$total.Range("A2").Copy()
$total.Range("A3").PasteSpecial(-4122)

# Insert the new 2022-Q3 totals row in row 2.
$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 14
$total.Range("D2").Value = 7.62
